# Update biosampleNumber (column C) for rows 2-9 so the values are
# consistent with rnaSampleNumber (column F): 1-8 -> 9-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 11
$ws.Range("C5").Value = 12
$ws.Range("C6").Value = 13
$ws.Range("C7").Value = 14
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 16

# Match the saved cursor position recorded in the workbook after editing.
$ws.Range("C10").Select() | Out-Null
